# cfb_weather.xlsx update
# Commit: Update cfb_weather.xlsx with Timestamp 2025-09-11T10:01:54.088523
#
# This script updates:
#  - the "FBS" sheet: refreshed odds/weather figures for several games and a
#    new run Timestamp (column AK) for every data row
#  - the "Other" sheet: refreshed odds/weather figures for several games,
#    plus the "Portland State vs Hawaii" game record moving up to directly
#    follow "Villanova vs Penn State" (rows 28-31 shift down by one, with the
#    Portland State/Hawaii data now occupying row 28)

$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

$newTimestamp = "2025-09-11T10:01:54.088523"

# ---------------------------------------------------------------------------
# FBS sheet: individual odds / weather figure corrections
# ---------------------------------------------------------------------------
$wsFBS.Range("Z2").Value = -114
$wsFBS.Range("Z3").Value = -110
$wsFBS.Range("Q5").Value = "E"
$wsFBS.Range("R6").Value = 0.1
$wsFBS.Range("Z7").Value = -105
$wsFBS.Range("Q8").Value = "ESE"
$wsFBS.Range("Z8").Value = -105
$wsFBS.Range("Z9").Value = -112
$wsFBS.Range("Z10").Value = -115
$wsFBS.Range("Q14").Value = "E"
$wsFBS.Range("R16").Value = 1.6
$wsFBS.Range("R21").Value = 0
$wsFBS.Range("Q22").Value = "WNW"
$wsFBS.Range("R24").Value = 1.3
$wsFBS.Range("Z25").Value = -105
$wsFBS.Range("Z27").Value = -115
$wsFBS.Range("Z28").Value = -110
$wsFBS.Range("AE31").Value = 0
$wsFBS.Range("Y31").Value = 42.5
$wsFBS.Range("Z31").Value = -105
$wsFBS.Range("R32").Value = 0
$wsFBS.Range("Q35").Value = "WNW"
$wsFBS.Range("Q36").Value = "E"
$wsFBS.Range("Z37").Value = -105
$wsFBS.Range("Q38").Value = "WNW"
$wsFBS.Range("R39").Value = 0.3
$wsFBS.Range("AE40").Value = 0.01801801801801802
$wsFBS.Range("Q40").Value = "WNW"
$wsFBS.Range("R40").Value = 0
$wsFBS.Range("Y40").Value = 56.5
$wsFBS.Range("Q41").Value = "WNW"
$wsFBS.Range("R42").Value = 0
$wsFBS.Range("AE43").Value = 0
$wsFBS.Range("Q43").Value = "NW"
$wsFBS.Range("Y43").Value = 51.5
$wsFBS.Range("Z43").Value = -115
$wsFBS.Range("Q45").Value = "E"

# Refresh the run Timestamp (column AK) on every data row (2-45)
for ($r = 2; $r -le 45; $r++) {
    $wsFBS.Range("AK$r").Value = $newTimestamp
}

# ---------------------------------------------------------------------------
# Other sheet: individual odds / weather figure corrections
# ---------------------------------------------------------------------------
$wsOther.Range("T11").Value = 1.1
$wsOther.Range("T12").Value = 0.1
$wsOther.Range("S21").Value = "NW"
$wsOther.Range("T34").Value = 0
$wsOther.Range("S35").Value = "E"
$wsOther.Range("T35").Value = 0
$wsOther.Range("S39").Value = "WNW"
$wsOther.Range("S41").Value = "ESE"
$wsOther.Range("S43").Value = "E"
$wsOther.Range("S49").Value = "WNW"

# ---------------------------------------------------------------------------
# Other sheet: rows 28-31 - "Portland State vs Hawaii" (previously row 31)
# now leads this block of games, so every row's data shifts down by one
# (row 31 -> row 28, row 28 -> row 29, row 29 -> row 30, row 30 -> row 31).
# Values are written explicitly per target row/column.
# ---------------------------------------------------------------------------

# Row 28 <- old row 31 (Portland State vs Hawaii)
$wsOther.Range("A28").Value = "Portland State vs Hawaii"
$wsOther.Range("B28").Value = "Hawaii"
$wsOther.Range("C28").Value = "Portland State"
$wsOther.Range("D28").Value = "SUN 09/14"
$wsOther.Range("E28").Value = "06:00 PM"
$wsOther.Range("F28").Value = "High"
$wsOther.Range("G28").Value = "NW-SE"
$wsOther.Range("H28").Value = "High"
$wsOther.Range("J28").Value = -45.40000000000001
$wsOther.Range("K28").Value = 75.33
$wsOther.Range("L28").Value = 53.03
$wsOther.Range("M28").Value = 14
$wsOther.Range("N28").Value = 2014
$wsOther.Range("O28").Value = "WSW"
$wsOther.Range("P28").Value = "WSW"
$wsOther.Range("Q28").Value = 80.60000000000001
$wsOther.Range("R28").Value = 9.6
$wsOther.Range("S28").Value = "WSW"
$wsOther.Range("U28").Value = -0.08
$wsOther.Range("V28").Value = -0.08
$wsOther.Range("W28").Value = -4.4
$wsOther.Range("X28").Value = "21.294294, -157.819338"

# Row 29 <- old row 28 (Davidson vs Tennessee Tech)
$wsOther.Range("A29").Value = "Davidson vs Tennessee Tech"
$wsOther.Range("B29").Value = "Tennessee Tech"
$wsOther.Range("C29").Value = "Davidson"
$wsOther.Range("D29").Value = "SAT 09/13"
$wsOther.Range("E29").Value = "12:00 PM"
$wsOther.Range("F29").Value = "Mid"
$wsOther.Range("G29").ClearContents()
$wsOther.Range("H29").ClearContents()
$wsOther.Range("J29").Value = 91.39511110000004
$wsOther.Range("K29").Value = 58.66
$wsOther.Range("L29").Value = 61.16
$wsOther.Range("M29").ClearContents()
$wsOther.Range("N29").Value = 1966
$wsOther.Range("O29").Value = "NE"
$wsOther.Range("P29").Value = "NE"
$wsOther.Range("Q29").Value = 80.48
$wsOther.Range("R29").Value = 4.7
$wsOther.Range("S29").Value = "NE"
$wsOther.Range("U29").Value = -0.06
$wsOther.Range("V29").Value = 0
$wsOther.Range("W29").ClearContents()
$wsOther.Range("X29").Value = "36.1780555, -85.506183"

# Row 30 <- old row 29 (Monmouth vs Charlotte)
$wsOther.Range("A30").Value = "Monmouth vs Charlotte"
$wsOther.Range("B30").Value = "Charlotte"
$wsOther.Range("C30").Value = "Monmouth"
$wsOther.Range("D30").Value = "SAT 09/13"
$wsOther.Range("E30").Value = "06:00 PM"
$wsOther.Range("F30").Value = "High"
$wsOther.Range("G30").Value = "NW-SE"
$wsOther.Range("H30").Value = "High"
$wsOther.Range("I30").ClearContents()
$wsOther.Range("J30").Value = 180.465509
$wsOther.Range("K30").Value = 62.01
$wsOther.Range("L30").Value = 55.89
$wsOther.Range("M30").Value = 5.4
$wsOther.Range("N30").Value = 2013
$wsOther.Range("O30").Value = "WSW"
$wsOther.Range("P30").Value = "WNW"
$wsOther.Range("Q30").Value = 80.48
$wsOther.Range("R30").Value = 5
$wsOther.Range("S30").Value = "WNW"
$wsOther.Range("T30").Value = 0
$wsOther.Range("U30").Value = -0.06
$wsOther.Range("V30").Value = 0
$wsOther.Range("W30").Value = -0.4
$wsOther.Range("X30").Value = "35.3105033, -80.7401194"

# Row 31 <- old row 30 (Indiana State vs Indiana)
$wsOther.Range("A31").Value = "Indiana State vs Indiana"
$wsOther.Range("B31").Value = "Indiana"
$wsOther.Range("C31").Value = "Indiana State"
$wsOther.Range("D31").Value = "FRI 09/12"
$wsOther.Range("E31").Value = "06:30 PM"
$wsOther.Range("F31").Value = "Mid"
$wsOther.Range("G31").Value = "N-S"
$wsOther.Range("H31").Value = "Med"
$wsOther.Range("I31").Value = "E"
$wsOther.Range("J31").Value = 76.12876889999998
$wsOther.Range("K31").Value = 54.97
$wsOther.Range("L31").Value = 55.49
$wsOther.Range("M31").Value = 6.2
$wsOther.Range("N31").Value = 1960
$wsOther.Range("O31").Value = "ESE"
$wsOther.Range("P31").Value = "SSE"
$wsOther.Range("Q31").Value = 80.23999999999999
$wsOther.Range("R31").Value = 4.1
$wsOther.Range("S31").Value = "ESE"
$wsOther.Range("T31").Value = 0
$wsOther.Range("U31").Value = -0.03
$wsOther.Range("V31").Value = 0
$wsOther.Range("W31").Value = -2.1
$wsOther.Range("X31").Value = "39.1808959, -86.5256217"
